$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column C keeps its text (quarter-date label) data type instead of
# being auto-recognised as a date literal by Excel's smart entry.
$ws.Range("C2:C64").NumberFormat = "@"

$rows = @(
    @{Row=2; C="01/10/2018"; D=270206}
    @{Row=3; C="01/01/2019"; D=268879}
    @{Row=4; C="01/04/2019"; D=269972}
    @{Row=5; C="01/07/2019"; D=272366}
    @{Row=6; C="01/10/2019"; D=277373}
    @{Row=7; C="01/01/2020"; D=273714}
    @{Row=8; C="01/04/2020"; D=257484}
    @{Row=9; C="01/07/2020"; D=257996}
    @{Row=10; C="01/10/2020"; D=259187}
    @{Row=11; C="01/01/2021"; D=256687}
    @{Row=12; C="01/04/2021"; D=256240}
    @{Row=13; C="01/07/2021"; D=255979}
    @{Row=14; C="01/10/2021"; D=254518}
    @{Row=15; C="01/01/2022"; D=257104}
    @{Row=16; C="01/04/2022"; D=268481}
    @{Row=17; C="01/07/2022"; D=281340}
    @{Row=18; C="01/10/2022"; D=287198}
    @{Row=19; C="01/01/2023"; D=284840}
    @{Row=20; C="01/04/2023"; D=287698}
    @{Row=21; C="01/07/2023"; D=295474}
    @{Row=22; C="01/10/2023"; D=301602}
    @{Row=23; C="01/10/2018"; D=42602}
    @{Row=24; C="01/01/2019"; D=42108}
    @{Row=25; C="01/04/2019"; D=42310}
    @{Row=26; C="01/07/2019"; D=42111}
    @{Row=27; C="01/10/2019"; D=43224}
    @{Row=28; C="01/01/2020"; D=42295}
    @{Row=29; C="01/04/2020"; D=38651}
    @{Row=30; C="01/07/2020"; D=37339}
    @{Row=31; C="01/10/2020"; D=38877}
    @{Row=32; C="01/01/2021"; D=37768}
    @{Row=33; C="01/04/2021"; D=38145}
    @{Row=34; C="01/07/2021"; D=38872}
    @{Row=35; C="01/10/2021"; D=38637}
    @{Row=36; C="01/01/2022"; D=38163}
    @{Row=37; C="01/04/2022"; D=40075}
    @{Row=38; C="01/07/2022"; D=42414}
    @{Row=39; C="01/10/2022"; D=43088}
    @{Row=40; C="01/01/2023"; D=43042}
    @{Row=41; C="01/04/2023"; D=43165}
    @{Row=42; C="01/07/2023"; D=44307}
    @{Row=43; C="01/10/2023"; D=45144}
    @{Row=44; C="01/10/2018"; D=1844}
    @{Row=45; C="01/01/2019"; D=1822}
    @{Row=46; C="01/04/2019"; D=1838}
    @{Row=47; C="01/07/2019"; D=1818}
    @{Row=48; C="01/10/2019"; D=1849}
    @{Row=49; C="01/01/2020"; D=1884}
    @{Row=50; C="01/04/2020"; D=1774}
    @{Row=51; C="01/07/2020"; D=1644}
    @{Row=52; C="01/10/2020"; D=1929}
    @{Row=53; C="01/01/2021"; D=1740}
    @{Row=54; C="01/04/2021"; D=1877}
    @{Row=55; C="01/07/2021"; D=1850}
    @{Row=56; C="01/10/2021"; D=1943}
    @{Row=57; C="01/01/2022"; D=1808}
    @{Row=58; C="01/04/2022"; D=1815}
    @{Row=59; C="01/07/2022"; D=1875}
    @{Row=60; C="01/10/2022"; D=1959}
    @{Row=61; C="01/01/2023"; D=1903}
    @{Row=62; C="01/04/2023"; D=1983}
    @{Row=63; C="01/07/2023"; D=1946}
    @{Row=64; C="01/10/2023"; D=1891}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# Re-apply the workbook's original (default) cell style so the text
# formatting trick above does not leave a lingering style change.
$ws.Range("C2:C64").Style = "Normal"

Write-Output "Updated $($rows.Count) rows in Sheet1 (columns C and D)."
